# Swap the "CPSC 4135" / "CPSC 4000" rows (A8:B8 <-> A9:B9) on the active sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "CPSC 4000"
$ws.Range("B8").Value = 0

$ws.Range("A9").Value = "CPSC 4135"
$ws.Range("B9").Value = 3
